$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-43 (Generation 0-41) -> Fitness 7618
$ws.Range("C2:C43").Value = 7618

# Rows 44-252 (Generation 42-250) -> Fitness 7534
$ws.Range("C44:C252").Value = 7534
